$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 (A2,B2) and row 4 (A4,B4)
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2

$ws.Range("A2").Value2 = $a4
$ws.Range("B2").Value2 = $b4
$ws.Range("A4").Value2 = $a2
$ws.Range("B4").Value2 = $b2
